# Applies the diff:
#  - rows 3,4,5,6,7,9 have their data permuted (the observation rows were
#    re-ordered / re-numbered) - only the cells whose value actually
#    changes are written.
#  - three brand-new observation rows (10, 11, 12) are appended.
#  - the sheet's used range/dimension grows from A1:AY9 to A1:AY12
#    automatically as a result of writing into row 12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- helpers -------------------------------------------------------

function Set-NumCell($row, $col, $val) {
    $ws.Cells.Item($row, $col).Value2 = $val
}

function Set-BoolCell($row, $col, $val) {
    $ws.Cells.Item($row, $col).Value2 = $val
}

# Plain text writer - safe for strings Excel would not mis-parse as a
# number/date/bool (the overwhelming majority of the text fields here).
function Set-StrCell($row, $col, $val) {
    $ws.Cells.Item($row, $col).Value2 = $val
}

# Text writer for strings that Excel's auto-detection would otherwise
# coerce into a number or a date (e.g. "25" or "2023-09-07"). Using a
# leading apostrophe (exactly like typing it into Excel) forces the
# cell to stay text.
function Set-ForcedStrCell($row, $col, $val) {
    $ws.Cells.Item($row, $col).Formula = "'" + $val
}

# ---- 1. permute rows 3,4,5,6,7,9 ------------------------------------
# new row 3 <- old row 9's data
Set-NumCell 3 1 111747186
Set-StrCell 3 16 "Lilla gruvan (Lilla gruvan), Ög"
Set-NumCell 3 17 575435.6246570286
Set-NumCell 3 18 6509856.898648335
Set-NumCell 3 19 2

# new row 4 <- old row 5's data
Set-NumCell 4 1 111749006
Set-NumCell 4 2 8377
Set-StrCell 4 4 "LC"
Set-NumCell 4 5 106545
Set-StrCell 4 6 "Mindre märgborre"
Set-StrCell 4 7 "Tomicus minor"
Set-StrCell 4 8 "(Hartig, 1834)"
Set-NumCell 4 17 575512.2089522779
Set-NumCell 4 18 6509825.662577543
Set-NumCell 4 19 2

# new row 5 <- old row 6's data
Set-NumCell 5 1 111747705
Set-NumCell 5 2 93067
Set-NumCell 5 5 2810
Set-StrCell 5 6 "Västlig hakmossa"
Set-StrCell 5 7 "Rhytidiadelphus loreus"
Set-StrCell 5 8 "(Hedw.) Warnst."
Set-NumCell 5 17 575459.4222356658
Set-NumCell 5 18 6509864.113963567

# new row 6 <- old row 7's data
Set-NumCell 6 1 111749860
Set-NumCell 6 2 78107
Set-StrCell 6 4 "NT"
Set-NumCell 6 5 6453
Set-StrCell 6 6 "Vedskivlav"
Set-StrCell 6 7 "Hertelidea botryosa"
Set-StrCell 6 8 "(Fr.) Printzen & Kantvilas"
Set-StrCell 6 16 "Älgsjöhåll (Älgsjöhåll), Ög"
Set-NumCell 6 17 575356.6078101217
Set-NumCell 6 18 6509772.251964441
Set-NumCell 6 19 1

# new row 7 <- old row 3's data
Set-NumCell 7 1 111749897
Set-NumCell 7 17 575336.6687912485
Set-NumCell 7 18 6509780.695668718

# new row 9 <- old row 4's data
Set-NumCell 9 1 111749343
Set-NumCell 9 17 575415.2450877089
Set-NumCell 9 18 6509807.674603676
Set-NumCell 9 19 1

# ---- 2. append new rows 10, 11, 12 ----------------------------------

# row 10
Set-NumCell 10 1 111964550
Set-NumCell 10 2 103288
Set-StrCell 10 3 "Ovaliderad"
Set-StrCell 10 4 "LC"
Set-NumCell 10 5 221144
Set-StrCell 10 6 "Grönpyrola"
Set-StrCell 10 7 "Pyrola chlorantha"
Set-StrCell 10 8 "Sw."
Set-ForcedStrCell 10 9 "25"
Set-StrCell 10 10 "plantor/tuvor"
Set-StrCell 10 11 "överblommad"
Set-StrCell 10 16 "Stenstorp SSO 1470 m, Ög"
Set-NumCell 10 17 575345.716659593
Set-NumCell 10 18 6509958.000975758
Set-NumCell 10 19 10
Set-StrCell 10 20 "Östergötland"
Set-StrCell 10 21 "Norrköping"
Set-StrCell 10 22 "Östergötland"
Set-StrCell 10 23 "Krokek"
Set-ForcedStrCell 10 25 "2023-09-07"
Set-StrCell 10 26 "00:00"
Set-ForcedStrCell 10 27 "2023-09-07"
Set-StrCell 10 28 "00:00"
Set-BoolCell 10 30 $false
Set-BoolCell 10 31 $false
Set-BoolCell 10 33 $false
Set-StrCell 10 35 "Äldre barrskog"
Set-StrCell 10 49 "Mirjam Ideström"
Set-StrCell 10 50 "Mirjam Ideström"

# row 11
Set-NumCell 11 1 111964621
Set-NumCell 11 2 93388
Set-StrCell 11 3 "Ovaliderad"
Set-StrCell 11 4 "LC"
Set-NumCell 11 5 2180
Set-StrCell 11 6 "Blåmossa"
Set-StrCell 11 7 "Leucobryum glaucum"
Set-StrCell 11 8 "(Hedw.) Ångstr."
Set-StrCell 11 16 "Stenstorp SSO 1660 m, Ög"
Set-NumCell 11 17 575609.0158921016
Set-NumCell 11 18 6509824.949736473
Set-NumCell 11 19 10
Set-StrCell 11 20 "Östergötland"
Set-StrCell 11 21 "Norrköping"
Set-StrCell 11 22 "Östergötland"
Set-StrCell 11 23 "Krokek"
Set-ForcedStrCell 11 25 "2023-09-07"
Set-StrCell 11 26 "00:00"
Set-ForcedStrCell 11 27 "2023-09-07"
Set-StrCell 11 28 "00:00"
Set-BoolCell 11 30 $false
Set-BoolCell 11 31 $false
Set-BoolCell 11 33 $false
Set-StrCell 11 35 "Barrskog"
Set-StrCell 11 49 "Mirjam Ideström"
Set-StrCell 11 50 "Mirjam Ideström"

# row 12
Set-NumCell 12 1 111964494
Set-NumCell 12 2 56414
Set-StrCell 12 3 "Ovaliderad"
Set-StrCell 12 4 "NT"
Set-NumCell 12 5 100049
Set-StrCell 12 6 "Spillkråka"
Set-StrCell 12 7 "Dryocopus martius"
Set-StrCell 12 8 "(Linnaeus, 1758)"
Set-StrCell 12 11 "adult"
Set-StrCell 12 13 "förbiflygande"
Set-StrCell 12 16 "Stenstorp SSO 1470 m, Ög"
Set-NumCell 12 17 575345.716659593
Set-NumCell 12 18 6509958.000975758
Set-NumCell 12 19 10
Set-StrCell 12 20 "Östergötland"
Set-StrCell 12 21 "Norrköping"
Set-StrCell 12 22 "Östergötland"
Set-StrCell 12 23 "Krokek"
Set-ForcedStrCell 12 25 "2023-09-07"
Set-StrCell 12 26 "10:30"
Set-ForcedStrCell 12 27 "2023-09-07"
Set-StrCell 12 28 "10:30"
Set-BoolCell 12 30 $false
Set-BoolCell 12 31 $false
Set-BoolCell 12 33 $false
Set-StrCell 12 35 "Äldre barrskog"
Set-StrCell 12 49 "Mirjam Ideström"
Set-StrCell 12 50 "Mirjam Ideström"
